# Implementación del caso de uso "Cambiar alumno de grupo" (CU 11):
# la tarea correspondiente (fila 14 de la hoja "Casos de Uso") se marca
# como terminada ("hecho") y se registran las horas consumidas el día 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Estatus de la tarea "CU 11 - Cambiar alumno de grupo." -> "hecho"
$ws.Range("F14").Value = "hecho"

# Horas consumidas el Día 5 (columna T, "Cons.") para esa tarea
$ws.Range("T14").Value = 4

# Deja la celda F14 como seleccionada, reflejando dónde quedó el
# usuario tras registrar el cambio.
$ws.Range("F14").Select()
